$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.871.21'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '1.642.00'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.79'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5063'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2587'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06438'
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.51'
$ws.Range("E10").Value = '  +5.25%  '
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.270'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.644.41'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.870.41'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5631'
$ws.Range("E15").Value = '  +2.28%  '
$ws.Range("D16").Value = '0.0₅7702'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.38'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '25.912.22'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.10'
$ws.Range("E20").Value = '  -1.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.373'
$ws.Range("E21").Value = '  -0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.929'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.122'
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.808'
$ws.Range("E25").Value = '  -5.79%  '
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1235'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.795'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.244'
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04956'
$ws.Range("E31").Value = '  +1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.294'
$ws.Range("E32").Value = '  +1.33%  '
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.573'
$ws.Range("E34").Value = '  +1.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9058'
$ws.Range("E36").Value = '  +0.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5550'
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = '1.133.20'
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.553'
$ws.Range("E39").Value = '  +0.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01566'
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9982'
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.510'
$ws.Range("E42").Value = '  -1.90%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8034'
$ws.Range("E43").Value = '  +0.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.07'
$ws.Range("E44").Value = '  +1.65%  '
$ws.Range("D45").Value = '1.781.82'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").Value = '0.0₈112'
$ws.Range("E46").Value = '  -5.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.73'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4291'
$ws.Range("E48").Value = '  -3.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.801'
$ws.Range("E49").Value = '  +2.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05044'
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  -0.21%  '
